# MB18_thresholds_5D.xlsx - "new config 5 Digit"
# Updates the threshold table on Sheet1 (rows 1-32, columns A:E) to the
# recalculated 5-digit configuration values and restores the saved
# selection (cell G6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    @(34, 33, 142, 143, 137),
    @(12380.95238095238, 5000, 5000, 7000, 10500),
    @(16190.476190476189, 14000, 15000, 7000, 18500),
    @(14761.904761904761, 13000, 14000, 7000, 19500),
    @(14761.904761904761, 14000, 15000, 7000, 19500),
    @(17142.857142857141, 19000, 16000, 7000, 21500),
    @(13333.333333333332, 12000, 10000, 7000, 15000),
    @(15238.095238095237, 17073.170731707316, 15000, 12195.121951219511, 19744.483159117302),
    @(15873.015873015869, 16666.666666666664, 17000, 11904.761904761903, 19274.376417233558),
    @(15503.875968992246, 16279.069767441859, 16279.069767441859, 11627.906976744185, 18826.13510520487),
    @(15151.515151515152, 15909.09090909091, 15909.09090909091, 11363.636363636364, 18398.268398268399),
    @(14814.814814814814, 15000, 15555.555555555555, 11111.111111111111, 18888.888888888887),
    @(15942.028985507248, 15217.391304347826, 15217.391304347826, 12000, 18478.260869565216),
    @(15602.836879432623, 14893.617021276596, 14893.617021276596, 14000, 18500),
    @(15238.095238095237, 14583.333333333332, 15500, 14000, 18500),
    @(12571.428571428572, 14285.714285714286, 14285.714285714286, 12000, 16500),
    @(14666.666666666668, 15000, 13000, 12000, 18000),
    @(14379.084967320263, 13725.490196078432, 13725.490196078432, 12000, 17000),
    @(12820.51282051282, 13461.538461538461, 13461.538461538461, 12000, 15000),
    @(12578.616352201256, 13207.54716981132, 13207.54716981132, 12500, 15000),
    @(13580.246913580246, 14000, 14000, 12000, 15500),
    @(13809.523809523809, 12727.272727272726, 12727.272727272726, 11000, 18500),
    @(13095.238095238095, 12499.999999999998, 12499.999999999998, 12500, 15178.571428571428),
    @(12865.497076023392, 12280.701754385964, 12280.701754385964, 11000, 13000),
    @(12643.67816091954, 12068.965517241379, 12068.965517241379, 12500, 13500),
    @(12429.378531073446, 11864.406779661016, 11864.406779661016, 12000, 14406.779661016948),
    @(13333.333333333332, 11666.666666666664, 11666.666666666664, 13000, 14166.666666666664),
    @(12021.857923497269, 11475.409836065573, 11475.409836065573, 9000, 15000),
    @(11827.956989247312, 11290.322580645161, 11290.322580645161, 9000, 13709.677419354839),
    @(11640.211640211641, 11111.111111111111, 11111.111111111111, 9000, 13492.063492063491),
    @(11458.333333333334, 10937.5, 10937.5, 9000, 13281.25),
    @(11282.051282051283, 10769.23076923077, 10769.23076923077, 9000, 13076.923076923076)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $values[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}

$ws.Range("G6").Select()
